$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4385.4736
$ws.Range("J64").Value = 4136.706
$ws.Range("L64").Value = 4136.706
$ws.Range("N64").Value = -4632.706

$ws.Range("H67").Value = 4385.4736
$ws.Range("J67").Value = 4136.706
$ws.Range("L67").Value = 4136.706
$ws.Range("N67").Value = -5852.706

$ws.Range("H76").Value = 3566.6667
$ws.Range("I76").Value = 3640
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 3640
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -3325
$ws.Range("N76").Value = -3830

$ws.Range("H79").Value = 3566.6667
$ws.Range("I79").Value = 3640
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 3640
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -2548
$ws.Range("N79").Value = -5384

$ws.Range("H111").Value = 2119.1738
$ws.Range("I111").Value = 1731.8889
$ws.Range("J111").Value = 2368.1428
$ws.Range("K111").Value = 5195.6667
$ws.Range("L111").Value = 7104.428400000001
$ws.Range("M111").Value = -2128.6667
$ws.Range("N111").Value = -13238.4284

$ws.Range("H137").Value = 2384333
$ws.Range("I137").Value = 5558032.5
$ws.Range("J137").Value = 4058.4583
$ws.Range("K137").Value = 16674097.5
$ws.Range("L137").Value = 12175.3749
$ws.Range("M137").Value = -16671547.5
$ws.Range("N137").Value = -17275.3749

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 45546560
$ws.Range("I61").Value = 52685350
$ws.Range("J61").Value = 334204.66
$ws.Range("K61").Value = 52685350
$ws.Range("L61").Value = 334204.66
$ws.Range("M61").Value = -52685138
$ws.Range("N61").Value = -334628.66

$ws.Range("H74").Value = 6212111.5
$ws.Range("I74").Value = 9288597
$ws.Range("J74").Value = 59139.668
$ws.Range("K74").Value = 9288597
$ws.Range("L74").Value = 59139.668
$ws.Range("M74").Value = -9287723
$ws.Range("N74").Value = -60887.668

$ws.Range("H77").Value = 6212111.5
$ws.Range("I77").Value = 9288597
$ws.Range("J77").Value = 59139.668
$ws.Range("K77").Value = 46442985
$ws.Range("L77").Value = 295698.34
$ws.Range("M77").Value = -46438617
$ws.Range("N77").Value = -304434.34

$ws.Range("H136").Value = 45546560
$ws.Range("I136").Value = 52685350
$ws.Range("J136").Value = 334204.66
$ws.Range("K136").Value = 158056050
$ws.Range("L136").Value = 1002613.98
$ws.Range("M136").Value = -158053500
$ws.Range("N136").Value = -1007713.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9733.821
$ws.Range("I86").Value = 16213.2
$ws.Range("J86").Value = 2257.6155
$ws.Range("K86").Value = 16213.2
$ws.Range("L86").Value = 2257.6155
$ws.Range("M86").Value = -15090.2
$ws.Range("N86").Value = -4503.6155

$ws.Range("H89").Value = 9733.821
$ws.Range("I89").Value = 16213.2
$ws.Range("J89").Value = 2257.6155
$ws.Range("K89").Value = 81066
$ws.Range("L89").Value = 11288.0775
$ws.Range("M89").Value = -75450
$ws.Range("N89").Value = -22520.0775

$ws.Range("H96").Value = 10647.556
$ws.Range("I96").Value = 6546.857
$ws.Range("K96").Value = 6546.857
$ws.Range("M96").Value = -3800.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 230414.48
$ws.Range("I31").Value = 51881.7
$ws.Range("K31").Value = 51881.7
$ws.Range("M31").Value = -51586.7

$ws.Range("H34").Value = 230414.48
$ws.Range("I34").Value = 51881.7
$ws.Range("K34").Value = 51881.7
$ws.Range("M34").Value = -51679.7

$ws.Range("H86").Value = 2067.7778
$ws.Range("I86").Value = 1300.3334
$ws.Range("J86").Value = 3602.6667
$ws.Range("K86").Value = 1300.3334
$ws.Range("L86").Value = 3602.6667
$ws.Range("M86").Value = -177.3334
$ws.Range("N86").Value = -5848.6667

$ws.Range("H89").Value = 2067.7778
$ws.Range("I89").Value = 1300.3334
$ws.Range("J89").Value = 3602.6667
$ws.Range("K89").Value = 6501.666999999999
$ws.Range("L89").Value = 18013.3335
$ws.Range("M89").Value = -885.6669999999995
$ws.Range("N89").Value = -29245.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1113.65
$ws.Range("I14").Value = 1113.65
$ws.Range("K14").Value = 3340.95
$ws.Range("M14").Value = -3167.95

$ws.Range("H70").Value = 2307.7144
$ws.Range("J70").Value = 3334.6667
$ws.Range("L70").Value = 10004.0001
$ws.Range("N70").Value = -10634.0001

$ws.Range("H73").Value = 2307.7144
$ws.Range("J73").Value = 3334.6667
$ws.Range("L73").Value = 10004.0001
$ws.Range("N73").Value = -12188.0001

$ws.Range("H80").Value = 1411.2858
$ws.Range("I80").Value = 1087.7646
$ws.Range("K80").Value = 3263.2938
$ws.Range("M80").Value = -2327.2938

$ws.Range("H83").Value = 1411.2858
$ws.Range("I83").Value = 1087.7646
$ws.Range("K83").Value = 9789.8814
$ws.Range("M83").Value = -5109.8814

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 1980
$ws.Range("I55").Value = 1980
$ws.Range("K55").Value = 1980
$ws.Range("M55").Value = -1653

$ws.Range("H70").Value = 45164.32
$ws.Range("I70").Value = 63812.234
$ws.Range("K70").Value = 63812.234
$ws.Range("M70").Value = -63542.234

$ws.Range("H73").Value = 45164.32
$ws.Range("I73").Value = 63812.234
$ws.Range("K73").Value = 63812.234
$ws.Range("M73").Value = -62876.234

$ws.Range("H80").Value = 4020
$ws.Range("I80").Value = 2800
$ws.Range("J80").Value = 4113.846
$ws.Range("K80").Value = 2800
$ws.Range("L80").Value = 4113.846
$ws.Range("M80").Value = -1802
$ws.Range("N80").Value = -6109.846

$ws.Range("H83").Value = 4020
$ws.Range("I83").Value = 2800
$ws.Range("J83").Value = 4113.846
$ws.Range("K83").Value = 14000
$ws.Range("L83").Value = 20569.23
$ws.Range("M83").Value = -9008
$ws.Range("N83").Value = -30553.23

$ws.Range("H132").Value = 54757.58
$ws.Range("I132").Value = 30905.176
$ws.Range("J132").Value = 257503
$ws.Range("K132").Value = 92715.52799999999
$ws.Range("L132").Value = 772509
$ws.Range("M132").Value = -90185.52799999999
$ws.Range("N132").Value = -777569

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 29999
$ws.Range("J98").Value = 29999
$ws.Range("L98").Value = 29999
$ws.Range("N98").Value = -35989

$ws.Range("H132").Value = 34627.97
$ws.Range("I132").Value = 2266.7144
$ws.Range("J132").Value = 336666.34
$ws.Range("K132").Value = 6800.1432
$ws.Range("L132").Value = 1009999.02
$ws.Range("M132").Value = -4270.1432
$ws.Range("N132").Value = -1015059.02

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 423.23334
$ws.Range("I107").Value = 353.82352
$ws.Range("J107").Value = 514
$ws.Range("K107").Value = 1061.47056
$ws.Range("L107").Value = 1542
$ws.Range("M107").Value = 858.52944
$ws.Range("N107").Value = -5382

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H132").Value = 55945.324
$ws.Range("I132").Value = 72777.42999999999
$ws.Range("J132").Value = 45699.695
$ws.Range("K132").Value = 218332.29
$ws.Range("L132").Value = 137099.085
$ws.Range("M132").Value = -215802.29
$ws.Range("N132").Value = -142159.085
